$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that become empty in the fixed forecaster output
$clearRefs = @(
    "C5", "E6", "G7", "I8", "K9", "M10", "O11", "Q12", "R13", "S13", "T14", "U14", "V14", "V15", "W15",
    "X15", "Y15", "Z15", "Y16", "Z16", "AA16", "AB16", "AC16", "AD16", "AC17", "AD17", "AE17", "AF17",
    "AG17", "AG18", "AH18", "AI18", "AJ18", "AK18", "AK19", "AL19", "AM19", "AN19", "AO19", "AO20",
    "AP20", "AQ20", "AR20", "AS20", "AS21", "AT21", "AU21", "AV21", "AW21", "AW22", "AX22", "AY22",
    "AZ22"
)
foreach ($ref in $clearRefs) { $ws.Range($ref).ClearContents() }

# Set corrected forecaster values
$setMap = @{
    "B1" = 39583
    "C1" = 39765
    "D1" = 39948
    "E1" = 40130
    "F1" = 40310
    "G1" = 40494
    "H1" = 40676
    "I1" = 40862
    "J1" = 41044
    "K1" = 41228
    "L1" = 41409
    "M1" = 41592
    "N1" = 41774
    "O1" = 41957
    "P1" = 42137
    "Q1" = 42321
    "R1" = 42503
    "S1" = 42689
    "T1" = 42867
    "U1" = 43053
    "V1" = 43145
    "W1" = 43235
    "X1" = 43326
    "Y1" = 43418
    "Z1" = 43510
    "AA1" = 43600
    "AB1" = 43691
    "AC1" = 43783
    "AD1" = 43875
    "AE1" = 43966
    "AF1" = 44068
    "AG1" = 44159
    "AH1" = 44251
    "AI1" = 44341
    "AJ1" = 44432
    "AK1" = 44525
    "AL1" = 44617
    "AM1" = 44706
    "AN1" = 44798
    "AO1" = 44890
    "AP1" = 44981
    "AQ1" = 45071
    "AR1" = 45163
    "AS1" = 45254
    "AT1" = 45345
    "AU1" = 45436
    "AV1" = 45534
    "AW1" = 45618
    "AX1" = 45713
    "AY1" = 45800
    "AZ1" = 45891
    "B3" = 1.003756253906252
    "C3" = 0.8212989654785341
    "D3" = 1.287693099940079
    "B4" = 1.003756253906229
    "C4" = 0.8527132153202777
    "D4" = 1.224010362214401
    "E4" = 1.183007486132071
    "F4" = 0.5167526861706184
    "D5" = 1.171780765338659
    "E5" = 1.129754998932886
    "F5" = 0.9718821796794952
    "G5" = 1.015842920196763
    "H5" = 1.187829657075357
    "F6" = 0.9030420588129306
    "G6" = 0.9131100874974818
    "H6" = 1.00065194548169
    "I6" = 0.9092565586104273
    "J6" = 1.071158385438342
    "H7" = 1.013831818446054
    "I7" = 0.8527132153202999
    "J7" = 0.912403143334517
    "K7" = 1.236730309040235
    "L7" = 0.9553801317191413
    "J8" = 0.8929747835070723
    "K8" = 1.175143392401168
    "L8" = 1.066801818459595
    "M8" = 1.029015928490629
    "N8" = 1.190496724073231
    "L9" = 1.084381878514384
    "M9" = 1.044063034969711
    "N9" = 1.154811676806311
    "O9" = 1.358148715145191
    "P9" = 1.5464392869869
    "N10" = 1.124712786946547
    "O10" = 1.200971162216891
    "P10" = 1.247870081683522
    "Q10" = 1.528208222695326
    "R10" = 1.701952652941463
    "P11" = 1.237299829452354
    "Q11" = 1.459520918994883
    "R11" = 1.637918813512695
    "S11" = 1.634928000057778
    "T11" = 1.580693894992691
    "R12" = 1.619591180870428
    "S12" = 1.639835143158264
    "T12" = 1.610567777412109
    "U12" = 1.67176973076042
    "V12" = 1.651937828695615
    "W12" = 1.646565058924154
    "X12" = 1.642460763882414
    "T13" = 1.624044593087226
    "U13" = 1.658462259692217
    "V13" = 1.63821551487775
    "W13" = 1.636439239090515
    "X13" = 1.62621273827539
    "Y13" = 1.603287858019664
    "Z13" = 1.504616869537312
    "AA13" = 1.619750436871126
    "AB13" = 1.242963308065193
    "W14" = 1.63821551487775
    "X14" = 1.630018843167003
    "Y14" = 1.599505522959754
    "Z14" = 1.599505522959732
    "AA14" = 1.669486277487398
    "AB14" = 1.082447181878954
    "AC14" = 0.8408455317168162
    "AD14" = 0.6216637650511503
    "AE14" = 0.2954364073068261
    "AF14" = -3.662861831460751
    "AA15" = 1.635941430694832
    "AB15" = 1.352809188882431
    "AC15" = 1.251834034633581
    "AD15" = 1.126729649114599
    "AE15" = 0.8326407735962826
    "AF15" = -1.110565553434917
    "AG15" = -1.875058665585216
    "AH15" = -3.604628722764358
    "AI15" = -1.564297238929013
    "AJ15" = 0.1010915562932313
    "AE16" = 0.9151237836431569
    "AF16" = -0.8254326060611827
    "AG16" = -1.301598234879375
    "AH16" = -2.856219939917704
    "AI16" = 0.1548119563699935
    "AJ16" = 6.64637963114707
    "AK16" = 5.03478667886097
    "AL16" = 5.220550987750228
    "AM16" = 5.937304773291885
    "AN16" = 5.793673192389748
    "AH17" = -2.405982383662564
    "AI17" = -1.150651722479712
    "AJ17" = 1.193655840614771
    "AK17" = 0.8711813995756934
    "AL17" = 1.043506288584606
    "AM17" = 2.150399152794202
    "AN17" = 1.728278600643907
    "AO17" = 2.399708479013141
    "AP17" = 0.156542203858212
    "AQ17" = -0.2621830498131694
    "AR17" = -0.3788601787194756
    "AL18" = 1.000460077742393
    "AM18" = 1.565839745803776
    "AN18" = 1.546011466763297
    "AO18" = 5.778276958487227
    "AP18" = 2.138492443986739
    "AQ18" = 1.878976297039481
    "AR18" = 2.181307424743695
    "AS18" = 0.8520283695166997
    "AT18" = 0.005756553697899847
    "AU18" = 0.0512320434504332
    "AV18" = 0.05771202657300911
    "AP19" = 2.483937396762181
    "AQ19" = 2.707799008349898
    "AR19" = 3.772966775860587
    "AS19" = 2.397640976910509
    "AT19" = -0.1096192596443557
    "AU19" = 0.2660756331863467
    "AV19" = 0.2051249733294291
    "AW19" = 0.299857156820571
    "AX19" = 0.4540776569412763
    "AY19" = 0.5998844096825495
    "AZ19" = 0.6062046309774693
    "AT20" = -0.1435426231004566
    "AU20" = 0.2967096123072999
    "AV20" = 0.1681072227950775
    "AW20" = -0.3196142230178167
    "AX20" = -0.1815195499670796
    "AY20" = 0.1733734969819434
    "AZ20" = 0.3877310837361314
    "AX21" = -0.2202168355973866
    "AY21" = 0.05961442020945285
    "AZ21" = 0.3721192462383671
}
foreach ($ref in $setMap.Keys) { $ws.Range($ref).Value = $setMap[$ref] }

# Remove trailing rows/column no longer present in the corrected output
$ws.Range("A23:A24").EntireRow.Delete()
$ws.Range("BA1").EntireColumn.Delete()